# Update the "Estado de Cuenta" worker rows (16-18).
# The previous record for YICETH COLON ANGULO (doc. 1143379752), which used to be
# the last row (18), is now the first row (16); the other two records shift down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1143379752"
$ws.Range("D16").Value = "YICETH COLON ANGULO"

$ws.Range("C17").Value = "9116598"
$ws.Range("D17").Value = "MARTIN EMILIO BERTEL GUZMAN"

$ws.Range("C18").Value = "1007640726"
$ws.Range("D18").Value = "SEIDER JOSE HERNANDEZ POLO"
